# Update countries & provincias Spain
# Applies the "Pais" sheet refresh: reorders four country-name pairs
# (which, because the sheet is sorted by "Casos totales" descending,
# also means the per-row statistics below them shift) and bumps the
# "Datos actualizados" timestamp, then writes the refreshed COVID
# figures for every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last refreshed" timestamp (row 1) ---------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 13:21"

# --- Per-row data (country label + the 7 stat columns B:H) -------------
# Row numbers refer to the fixed position in the worksheet; the country
# label in column A is written explicitly for every row so the four
# adjacent re-orderings (Rumania/Singapur, Republica de Chipre/Niger,
# Vietnam/Bahamas, Papua Nueva Guinea/Butan) land correctly alongside
# their refreshed numbers.
$updates = @(
    @{Row=13;  Country='Iran';                   B=309437; C=2685; D=268102; E=24145; F=0; G=208; H=17190}
    @{Row=27;  Country='Catar';                   B=111107; C=196;  D=107779; E=3151;  F=0; G=3;   H=177}
    @{Row=41;  Country='Kuwait';                  B=67911;  C=463;  D=59213;  E=8241;  F=0; G=4;   H=457}
    @{Row=43;  Country='Emiratos Arabes Unidos';  B=60999;  C=239;  D=54615;  E=6033;  F=0; G=0;   H=351}
    @{Row=45;  Country='Rumania';                 B=53186;  C=1075; D=27592;  E=23181; F=0; G=34;  H=2413}
    @{Row=46;  Country='Singapur';                B=52825;  C=313;  D=46740;  E=6058;  F=0; G=0;   H=27}
    @{Row=49;  Country='Polonia';                 B=46894;  C=548;  D=34709;  E=10454; F=0; G=10;  H=1731}
    @{Row=58;  Country='Suiza';                   B=35550;  C=138;  D=31300;  E=2269;  F=0; G=0;   H=1981}
    @{Row=68;  Country='Nepal';                   B=20332;  C=246;  D=14603;  E=5672;  F=0; G=1;   H=57}
    @{Row=83;  Country='Madagascar';              B=11528;  C=255;  D=8444;   E=2970;  F=0; G=7;   H=114}
    @{Row=85;  Country='Senegal';                 B=10344;  C=60;   D=6838;   E=3297;  F=0; G=0;   H=209}
    @{Row=121; Country='Mali';                    B=2535;   C=0;    D=1941;   E=470;   F=0; G=0;   H=124}
    @{Row=147; Country='Republica de Chipre';     B=1139;   C=15;   D=852;    E=268;   F=0; G=0;   H=19}
    @{Row=148; Country='Niger';                   B=1136;   C=0;    D=1028;   E=39;    F=0; G=0;   H=69}
    @{Row=154; Country='Malta';                   B=860;    C=15;   D=666;    E=185;   F=0; G=0;   H=9}
    @{Row=161; Country='Vietnam';                 B=620;    C=30;   D=373;    E=242;   F=0; G=2;   H=5}
    @{Row=162; Country='Bahamas';                 B=599;    C=0;    D=91;     E=494;   F=0; G=0;   H=14}
    @{Row=189; Country='Papua Nueva Guinea';      B=110;    C=19;   D=34;     E=74;    F=0; G=0;   H=2}
    @{Row=190; Country='Butan';                   B=102;    C=1;    D=89;     E=13;    F=0; G=0;   H=0}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.Country
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = $u.F
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
}
